# Attendance-scanner re-upload: the scan log row is updated with a real
# (numeric) scan time instead of a plain text time-stamp, so Excel can
# format/sort it as a time value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 held the literal text "14:13:09" (t="str"). Replace it with the actual
# Excel time serial for the new scan (11:03:15 -> 0.46059027777777778 of a
# 24h day) and format it as h:mm:ss (built-in numFmtId 21).
$ws.Range("D2").Value = 0.46059027777777778
$ws.Range("D2").NumberFormat = "h:mm:ss"
$ws.Range("D2").Font.Color = 0

# The re-saved workbook also carries a trailing, otherwise-empty row right
# below the data (same row height/font as the data rows) and leaves the
# selection sitting on the newly written time cell.
$ws.Rows.Item(3).RowHeight = 15.5
$ws.Rows.Item(3).Font.Size = 12
$ws.Range("D2").Select() | Out-Null
